$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string text: "Additional lipid extraction efficiency [%]" -> "Bagasse lipid extraction efficiency [%]" (used by B6)
$ws.Range("B6").Value = "Bagasse lipid extraction efficiency [%]"

# Update Spearman correlation values in rows 4-13 (columns C, D, E, F, H, I, J)

# Row 4
$ws.Range("C4").Value2 = 0.01444236105902648
$ws.Range("D4").Value2 = 0.9054916372909324
$ws.Range("E4").Value2 = -0.9773569339233482
$ws.Range("F4").Value2 = 0.8608070201755045
$ws.Range("H4").Value2 = 0.109868746718668
$ws.Range("I4").Value2 = 0.6153678841971051
$ws.Range("J4").Value2 = 0.9290407260181506

# Row 5
$ws.Range("C5").Value2 = 0.1106382659566489
$ws.Range("D5").Value2 = 0.2132198304957624
$ws.Range("E5").Value2 = -0.04777619440486013
$ws.Range("F5").Value2 = 0.1837830945773645
$ws.Range("H5").Value2 = 0.2848241206030152
$ws.Range("I5").Value2 = 0.1814535363384085
$ws.Range("J5").Value2 = 0.112870321758044

# Row 6
$ws.Range("C6").Value2 = 0.1157818945473637
$ws.Range("D6").Value2 = 0.3366729168229206
$ws.Range("E6").Value2 = -0.01833195829895748
$ws.Range("F6").Value2 = -0.1367209180229506
$ws.Range("H6").Value2 = 0.06141303532588316
$ws.Range("I6").Value2 = -0.1361899047476187
$ws.Range("J6").Value2 = 0.06338858471461788

# Row 7
$ws.Range("C7").Value2 = 0.1319132978324458
$ws.Range("D7").Value2 = 0.09170779269481738
$ws.Range("E7").Value2 = 0.2503607590189755
$ws.Range("F7").Value2 = 0.2926633165829146
$ws.Range("H7").Value2 = 0.8045796144903624
$ws.Range("I7").Value2 = 0.6976464411610291
$ws.Range("J7").Value2 = -0.07246681167029177

# Row 8
$ws.Range("C8").Value2 = 0.7028215705392636
$ws.Range("D8").Value2 = 0.04548113702842571
$ws.Range("E8").Value2 = -0.02687617190429761
$ws.Range("F8").Value2 = 0.07186379659491487
$ws.Range("H8").Value2 = 0.04589064726618167
$ws.Range("I8").Value2 = 0.06453761344033601
$ws.Range("J8").Value2 = 0.02938423460586515

# Row 9
$ws.Range("C9").Value2 = 0.4342368559213981
$ws.Range("D9").Value2 = -0.09235580889522239
$ws.Range("E9").Value2 = 0.08013350333758347
$ws.Range("F9").Value2 = -0.0853206330158254
$ws.Range("H9").Value2 = -0.01752793819845496
$ws.Range("I9").Value2 = -0.02825170629265732
$ws.Range("J9").Value2 = -0.08870321758043952

# Row 10
$ws.Range("C10").Value2 = 0.09866796669916748
$ws.Range("D10").Value2 = -0.002994074851871297
$ws.Range("E10").Value2 = -0.04183604590114753
$ws.Range("F10").Value2 = 0.01905497637440936
$ws.Range("H10").Value2 = -0.002796069901747544
$ws.Range("I10").Value2 = -0.01152028800720018
$ws.Range("J10").Value2 = 0.04182854571364285

# Row 11
$ws.Range("C11").Value2 = 0.2010365259131479
$ws.Range("D11").Value2 = 0.0634215855396385
$ws.Range("E11").Value2 = -0.03498987474686868
$ws.Range("F11").Value2 = 0.07256581414535364
$ws.Range("H11").Value2 = 0.08938723468086704
$ws.Range("I11").Value2 = 0.1057781444536114
$ws.Range("J11").Value2 = 0.08425260631515789

# Row 12
$ws.Range("C12").Value2 = 0.1235865896647416
$ws.Range("D12").Value2 = 0.08455411385284634
$ws.Range("E12").Value2 = 0.09969549238730969
$ws.Range("F12").Value2 = 0.1820535513387835
$ws.Range("H12").Value2 = 0.4177814445361135
$ws.Range("I12").Value2 = -0.07653491337283434
$ws.Range("J12").Value2 = -0.01932948323708093

# Row 13
$ws.Range("C13").Value2 = -0.3042901072526814
$ws.Range("D13").Value2 = 0.002371559288982225
$ws.Range("E13").Value2 = -0.003367584189604741
$ws.Range("F13").Value2 = 0.02346058651466287
$ws.Range("H13").Value2 = 0.02072601815045376
$ws.Range("I13").Value2 = 0.0548068701717543
$ws.Range("J13").Value2 = -0.01204980124503113
